$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-order / update the team member table ---
# Row 1: Phan Huu Phuoc (now marked as team leader) moves to the top
$ws.Range("B1").Value = 21200333
$ws.Range("C1").Value = "Phan Hữu Phước (Trưởng nhóm)"
$ws.Range("D1").Value = "Tab Đạo hàm"

# Row 2: Vo Thanh Danh
$ws.Range("B2").Value = 21200275
$ws.Range("C2").Value = "Võ Thành Danh"
$ws.Range("D2").Value = "Tab Nội suy"

# Row 3: Nguyen Tien Dat
$ws.Range("B3").Value = 21200278
$ws.Range("C3").Value = "Nguyễn Tiến Đạt"
$ws.Range("D3").Value = "Tab Tìm nghiệm"

# Row 5: Le Doan Phu Sang (group intro moved here, from the regression tab)
$ws.Range("B5").Value = 21200346
$ws.Range("C5").Value = "Lê Đoàn Phú Sang"
$ws.Range("D5").Value = "Tab Tìm nghiệm + Giới thiệu nhóm"

# Row 4: Tran Tuan Kiet (group intro removed from this tab)
$ws.Range("B4").Value = 21200304
$ws.Range("C4").Value = "Trần Tuấn Kiệt"
$ws.Range("D4").Value = "Tab Hồi quy"

# Row 6: Le Minh Thanh (unchanged)
$ws.Range("B6").Value = 21200351
$ws.Range("C6").Value = "Lê Minh Thành"
$ws.Range("D6").Value = "Tab Tích phân"

# --- Column widths ---
$ws.Range("C1").EntireColumn.ColumnWidth = 28.15
$ws.Range("D1").EntireColumn.ColumnWidth = 28.15

# --- Selection ---
$ws.Range("E5").Select()
